# Fix(joueur) : Placement du joueur
# Fill in the "Déplacement" day-block (row 26 header + row 27 task) of the
# Journal sheet that was previously left empty, and move the active
# selection to the next empty line (E27) as left by the author.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")
$ws.Activate()

# Row 26 : new day entry "Déplacment" (header line of the day block)
$ws.Range("A26").Value = "Déplacment"
$ws.Range("C26").Value = 50
$ws.Range("E26").Value = "Déplacement du joueur avec tous les problemes qui on suivit"

# Row 27 : task line "Explication"
$ws.Range("A27").Value = "Explication"
$ws.Range("C27").Value = 20
$ws.Range("E27").Value = "Explication des déplacement par Tony + celle du prof"

# Leave the selection on the next free cell, like in the authored workbook
$ws.Range("E27").Select()
